# Quarterly Financials update: insert a new (most-recent) quarter column
# before column D, shifting the existing D:K data to E:L, then populate
# the new column D with the latest quarter's figures. Also corrects a
# handful of historical "Capital Expenditures" (row 91) data points that
# were revised upstream.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1) Insert a new blank column before D - this shifts D:K -> E:L
#    (and updates the sheet dimension / row spans automatically).
$ws.Range("D1").EntireColumn.Insert()

# 2) Copy the number/cell formatting from the (now shifted) E column
#    into the new D column so the new cells match their row's style,
#    restricted to the three data blocks that actually contain data.
$ws.Range("E7:E35").Copy()
$ws.Range("D7:D35").PasteSpecial(-4122)

$ws.Range("E38:E77").Copy()
$ws.Range("D38:D77").PasteSpecial(-4122)

$ws.Range("E80:E102").Copy()
$ws.Range("D80:D102").PasteSpecial(-4122)

# 3) Populate the new column D (newest quarter, period ending 2018-06-30)
#    -- Income Statement block --
$ws.Range("D7").Value = 43281
$ws.Range("D8").Value = 3844900
$ws.Range("D9").Value = 2725300
$ws.Range("D10").Value = 1119600
$ws.Range("D12").Value = "NA"
$ws.Range("D13").Value = 0
$ws.Range("D14").Value = 0
$ws.Range("D15").Value = 0
$ws.Range("D17").Value = 3322000
$ws.Range("D18").Value = 522900
$ws.Range("D20").Value = 218100
$ws.Range("D21").Value = 950200
$ws.Range("D22").Value = 195400
$ws.Range("D23").Value = 545600
$ws.Range("D24").Value = 136700
$ws.Range("D25").Value = 0
$ws.Range("D26").Value = 408900
$ws.Range("D27").Value = 133800
$ws.Range("D28").Value = 0
$ws.Range("D29").Value = "NA"
$ws.Range("D30").Value = 0
$ws.Range("D31").Value = 0
$ws.Range("D32").Value = -218100
$ws.Range("D33").Value = 133800
$ws.Range("D34").Value = 0
$ws.Range("D35").Value = 133800

# -- Balance Sheet block --
$ws.Range("D38").Value = 43281
$ws.Range("D41").Value = 1856400
$ws.Range("D42").Value = 78500
$ws.Range("D43").Value = 1245100
$ws.Range("D44").Value = 828200
$ws.Range("D45").Value = 122200
$ws.Range("D46").Value = 4130400
$ws.Range("D47").Value = 5251700
$ws.Range("D48").Value = 5158600
$ws.Range("D49").Value = 4691300
$ws.Range("D50").Value = 0
$ws.Range("D51").Value = 0
$ws.Range("D52").Value = 645900
$ws.Range("D53").Value = 0
$ws.Range("D54").Value = 19877900
$ws.Range("D57").Value = 1388900
$ws.Range("D58").Value = 1832200
$ws.Range("D59").Value = 479500
$ws.Range("D60").Value = 3700600
$ws.Range("D61").Value = 5955800
$ws.Range("D62").Value = 1778600
$ws.Range("D63").Value = 0
$ws.Range("D64").Value = 0
$ws.Range("D65").Value = 0
$ws.Range("D66").Value = 16772900
$ws.Range("D68").Value = 0
$ws.Range("D69").Value = 0
$ws.Range("D70").Value = 0
$ws.Range("D71").Value = 0
$ws.Range("D72").Value = 3466300
$ws.Range("D73").Value = 0
$ws.Range("D74").Value = 0
$ws.Range("D75").Value = 0
$ws.Range("D76").Value = 3105000
$ws.Range("D77").Value = 0

# -- Cash Flow Statement block --
$ws.Range("D80").Value = 43281
$ws.Range("D81").Value = 133800
$ws.Range("D83").Value = 209200
$ws.Range("D84").Value = 0
$ws.Range("D85").Value = 0
$ws.Range("D86").Value = 0
$ws.Range("D87").Value = 0
$ws.Range("D88").Value = 0
$ws.Range("D89").Value = 398500
$ws.Range("D91").Value = -202100
$ws.Range("D92").Value = 0
$ws.Range("D93").Value = 0
$ws.Range("D94").Value = -325500
$ws.Range("D96").Value = -30400
$ws.Range("D97").Value = 0
$ws.Range("D98").Value = 0
$ws.Range("D99").Value = 0
$ws.Range("D100").Value = -118600
$ws.Range("D101").Value = -105900
$ws.Range("D102").Value = -151500

# 4) Row 91 ("Capital Expenditures") historical values were also revised
#    upstream - overwrite the full row with the corrected figures.
$ws.Range("D91").Value = -202100
$ws.Range("E91").Value = -294500
$ws.Range("F91").Value = -319600
$ws.Range("G91").Value = -129700
$ws.Range("H91").Value = -189800
$ws.Range("I91").Value = -39900
$ws.Range("J91").Value = -37300
$ws.Range("K91").Value = -309600
$ws.Range("L91").Value = -326800
